$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q6)
$ws.Range("B7").Value = -0.1062043836119984
$ws.Range("C7").Value = 3.214038737436144
$ws.Range("D7").Value = 28.42386334192099
$ws.Range("E7").Value = 5.331403505824802
$ws.Range("F7").Value = 5.401897041383672
$ws.Range("G7").Value = 38

# Row 8 (Q7)
$ws.Range("B8").Value = 0.1126760984253048
$ws.Range("C8").Value = 3.736719421784253
$ws.Range("D8").Value = 28.69179275025208
$ws.Range("E8").Value = 5.356472043262438
$ws.Range("F8").Value = 5.429156327292642
$ws.Range("G8").Value = 37

# Row 9 (Q8)
$ws.Range("B9").Value = 0.3025403896483665
$ws.Range("C9").Value = 4.219950309719327
$ws.Range("D9").Value = 42.57676554919534
$ws.Range("E9").Value = 6.525087397820457
$ws.Range("F9").Value = 6.68739858648746
$ws.Range("G9").Value = 20

# Row 10 (Q9)
$ws.Range("B10").Value = -1.273655624497557
$ws.Range("C10").Value = 4.102290063633903
$ws.Range("D10").Value = 42.28739103502502
$ws.Range("E10").Value = 6.502875597381901
$ws.Range("F10").Value = 6.637315603262485
$ws.Range("G10").Value = 13

# Row 11
$ws.Range("B11").Value = -3.449580058531092
$ws.Range("C11").Value = 5.366474568182582
$ws.Range("D11").Value = 52.70517324201703
$ws.Range("E11").Value = 7.259832865983695
$ws.Range("F11").Value = 7.141915942326126
$ws.Range("G11").Value = 5
